# Update the output of answers-of-two-digit_number_divided_by_one-digit_number.docx
# to reflect a freshly generated answer key (for date 2024-03-07 Thursday).

$d = $word.ActiveDocument

# 1) Update the date heading paragraph.
$d.Paragraphs.Item(1).Range.Text = "2024-03-07 Thursday"

# 2) Update each answer cell in the single table. The table has 20 rows but
#    only every 4th row (1, 5, 9, 13, 17) holds the 5 answer cells; the rows
#    in between are blank "working space" rows.
$t = $d.Tables.Item(1)

$answers = @(
  @("39÷9=4, 3", "26÷8=3, 2", "17÷3=5, 2", "49÷2=24, 1", "45÷8=5, 5"),
  @("36÷4=9, 0", "73÷8=9, 1", "40÷6=6, 4", "56÷7=8, 0", "42÷8=5, 2"),
  @("46÷4=11, 2", "17÷5=3, 2", "31÷8=3, 7", "51÷5=10, 1", "20÷8=2, 4"),
  @("43÷3=14, 1", "80÷9=8, 8", "68÷4=17, 0", "24÷3=8, 0", "87÷5=17, 2"),
  @("56÷8=7, 0", "49÷8=6, 1", "64÷4=16, 0", "76÷6=12, 4", "72÷9=8, 0")
)

$dataRows = @(1, 5, 9, 13, 17)

for ($i = 0; $i -lt $dataRows.Length; $i++) {
  $row = $dataRows[$i]
  $rowAnswers = $answers[$i]
  for ($c = 1; $c -le 5; $c++) {
    $t.Cell($row, $c).Range.Text = $rowAnswers[$c - 1]
  }
}
